# Weekly update: a new price observation for "Ciboulette" (Hortaliza) at
# "Vega Central Mapocho de Santiago" is inserted as row 633, pushing the
# existing rows 633:664 down to 634:665 (and extending the sheet's used
# range from R664 to R665).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 633:664 down by inserting a fresh row above the current row 633.
$ws.Rows.Item(633).Insert()

# Populate the newly inserted row 633 with the new observation.
$ws.Range("A633").Value = 9
$ws.Range("B633").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C633").Value = "Metropolitana"
$ws.Range("D633").Value = 45267
$ws.Range("E633").Value = 13
$ws.Range("F633").Value = 100112039
$ws.Range("G633").Value = "Ciboulette"
$ws.Range("H633").Value = "Sin especificar"
$ws.Range("I633").Value = "Primera"
$ws.Range("J633").Value = 430
$ws.Range("K633").Value = 1000
$ws.Range("L633").Value = 1200
$ws.Range("M633").Value = 1100
$ws.Range("N633").Value = "`$/docena de atados"
$ws.Range("O633").Value = "Región Metropolitana"
$ws.Range("P633").Value = 367
$ws.Range("Q633").Value = 3
$ws.Range("R633").Value = "Hortaliza"

# Match the date-serial style used by the rest of column D (row 634, the
# row that was previously at 633, already carries that style after the
# insert copied formatting down -- this call just makes it explicit/robust).
$ws.Range("D633").NumberFormat = $ws.Range("D634").NumberFormat
